$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '51.811.96'
$ws.Range("E2").Value = '  +0.09%  '

# Row 3
$ws.Range("D3").Value = '2.785.89'
$ws.Range("E3").Value = '  -1.12%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '357.24'
$ws.Range("E5").Value = '  +0.62%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '110.46'
$ws.Range("E6").Value = '  -0.86%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.559'
$ws.Range("E7").Value = '  -1.10%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.590'
$ws.Range("E9").Value = '  -1.20%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.04'
$ws.Range("E10").Value = '  -1.55%  '

# Row 11
$ws.Range("E11").Value = '  +2.25%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0849'
$ws.Range("E12").Value = '  -0.55%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.54'
$ws.Range("E13").Value = '  -1.73%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.63'
$ws.Range("E14").Value = '  -1.44%  '

# Row 15
$ws.Range("D15").Value = '3.231.25'
$ws.Range("E15").Value = '  -0.89%  '

# Row 16
$ws.Range("D16").Value = '2.787.34'
$ws.Range("E16").Value = '  -0.89%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.941'
$ws.Range("E17").Value = '  +2.51%  '

# Row 18
$ws.Range("D18").Value = '51.793.83'
$ws.Range("E18").Value = '  +0.15%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.47'
$ws.Range("E19").Value = '  -1.49%  '

# Row 20
$ws.Range("B20").Value = 'ImmutableX'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.06'
$ws.Range("E20").Value = '  -1.87%  '

# Row 21
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.33'
$ws.Range("E21").Value = '  -0.01%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0974'
$ws.Range("E22").Value = '  -1.65%  '

# Row 23
$ws.Range("E23").Value = '  +0.59%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '270.06'
$ws.Range("E24").Value = '  +0.85%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.76'
$ws.Range("E25").Value = '  -1.22%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.44'
$ws.Range("E26").Value = '  -1.63%  '

# Row 27
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.06%  '

# Row 28
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.166'
$ws.Range("E28").Value = '  +18.50%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.25'
$ws.Range("E29").Value = '  -0.41%  '

# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.23'
$ws.Range("E30").Value = '  -0.57%  '

# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.25'
$ws.Range("E31").Value = '  +6.69%  '

# Row 32
$ws.Range("B32").Value = 'OKB'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '52.06'
$ws.Range("E32").Value = '  -0.98%  '

# Row 33
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.94'
$ws.Range("E33").Value = '  +1.22%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0454'
$ws.Range("E34").Value = '  -7.56%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0845'
$ws.Range("E35").Value = '  +0.25%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.14'
$ws.Range("E36").Value = '  -4.25%  '

# Row 37
$ws.Range("E37").Value = '  +0.05%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.75'
$ws.Range("E38").Value = '  +2.63%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.18'
$ws.Range("E39").Value = '  -3.12%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.96'
$ws.Range("E40").Value = '  -3.34%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.58'
$ws.Range("E41").Value = '  +1.68%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.114'
$ws.Range("E42").Value = '  -1.84%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.23'
$ws.Range("E43").Value = '  -1.92%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.97'
$ws.Range("E44").Value = '  -4.13%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.78'
$ws.Range("E45").Value = '  -6.26%  '

# Row 46
$ws.Range("D46").Value = '2.084.92'
$ws.Range("E46").Value = '  -0.41%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.29'
$ws.Range("E47").Value = '  -1.28%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.28'
$ws.Range("E48").Value = '  +1.09%  '

# Row 49
$ws.Range("E49").Value = '  -4.12%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.940'
$ws.Range("E50").Value = '  -4.84%  '

# Row 51
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.64'
$ws.Range("E51").Value = '  -4.29%  '
